# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    6  = 3013
    8  = 2073
    10 = 106
    11 = 921
    14 = 229
    17 = 32
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
